$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF for Wins / Losses / Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season record values for every data row (2-50)
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 82
    $ws.Cells.Item($row, 31).Value = 80
    $ws.Cells.Item($row, 32).Value = 0
}
